# Apply edits to 'To Touch the Moon.xlsx' per the target diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the 'Links' sheet (1st tab) with additional link rows.
# ---------------------------------------------------------------
$wsLinks = $wb.Worksheets.Item("Links")
$wsLinks.Range("A2").Value = "Picture to ASCII Generator"
$wsLinks.Range("L2").Value = "Description"
$wsLinks.Range("A5").Value = "https://dwarffortresswiki.org/Tileset_repository"
$wsLinks.Range("L5").Value = "place to get fonts"
$wsLinks.Range("A6").Value = "https://www.text-image.com/convert/ascii.html"
$wsLinks.Range("L6").Value = "image to ascii converter"
$wsLinks.Range("O6").Value = "https://share.text-image.com/67e70ca32cf421d6"
$wsLinks.Range("A9").Value = "https://github.com/Thraka/SadConsole"
$wsLinks.Range("L9").Value = "SadConsole github"
$wsLinks.Range("A12").Value = "https://markjames.dev/2020-05-21-making-a-roguelike-in-c-with-gorogue-sadconsole-part-one/"
$wsLinks.Range("L12").Value = "Sad console tutorial"
$wsLinks.Range("A14").Value = "https://code2d.wordpress.com/sadconsole-tutorials/"
$wsLinks.Range("L14").Value = "Sad  Console tutorial"
$wsLinks.Range("A17").Value = "http://sadconsole.com/v9/api/"
$wsLinks.Range("L17").Value = "Link to SadConsole Documentation"

# ---------------------------------------------------------------
# 2. Add new weapon/item names to the 'Names' sheet.
# ---------------------------------------------------------------
$wsNames = $wb.Worksheets.Item("Names")
$wsNames.Range("I16").Value = "Blade of Pandemonium"
$wsNames.Range("I17").Value = "Scepter of the Heavenly King"
$wsNames.Range("I18").Value = "Roxy's Holy Relic"
$wsNames.Range("I19").Value = "Aylas Holy Scarf"

# ---------------------------------------------------------------
# 3. Create the new 'Monsters' worksheet (placed after 'Quests').
#    Creating it first means it is allocated sheetId 10.
# ---------------------------------------------------------------
$wsMonsters = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMonsters.Name = "Monsters"
$wsMonsters.Range("C1").Value = "Almost always neutral"
$wsMonsters.Range("C2").Value = "Machines made from past or future civilizations"
$wsMonsters.Range("C3").Value = "to angry to die"
$wsMonsters.Range("C4").Value = "Servants of Luna - Beings of Chaos"
$wsMonsters.Range("C5").Value = "Servants of Sol - Beings of Order"
$wsMonsters.Range("C6").Value = "have human like traits"
$wsMonsters.Range("C7").Value = "sub category of humanoid - tall beings"
$wsMonsters.Range("C8").Value = "Strong monsters"
$wsMonsters.Range("A11").Value = "Monster Prefix"
$wsMonsters.Range("Q11").Value = "Monster Suffix"
$wsMonsters.Range("A13").Value = "Raging"
$wsMonsters.Range("C13").Value = "Generated Hostile regardless of alignment - unable to be calmed down"
$wsMonsters.Range("Q13").Value = "Paladin"
$wsMonsters.Range("S13").Value = "will spawn friendly for neutral and lawful, 50/50 for chaos. Killing them will lose favor with that god unless chaotic and being chaotic. Luna finds it amusing when her minions fight each other."
$wsMonsters.Range("A14").Value = "Serene"
$wsMonsters.Range("C14").Value = "Generated peaceful regarless of alignment"
$wsMonsters.Range("Q14").Value = "of  *element type*"
$wsMonsters.Range("S14").Value = "extra 10% damage done is added to attack. I.E if monster roles 50 damage before calculations they will do +5 *elemental* damage added to final attack."
$wsMonsters.Range("A15").Value = "Master"
$wsMonsters.Range("C15").Value = " +5 Levels of expereince to monster"
$wsMonsters.Range("Q15").Value = " "
$wsMonsters.Range("A21").Value = "Animals"
$wsMonsters.Range("C21").Value = "Constructs"
$wsMonsters.Range("E21").Value = "Undead"
$wsMonsters.Range("G21").Value = "Demons"
$wsMonsters.Range("I21").Value = "Angels"
$wsMonsters.Range("K21").Value = "Humanoid"
$wsMonsters.Range("M21").Value = "Giants"
$wsMonsters.Range("O21").Value = "Dragons"
$wsMonsters.Range("Q21").Value = "Monsters"
$wsMonsters.Range("A23").Value = "Lion"
$wsMonsters.Range("C23").Value = "Iron Golem"
$wsMonsters.Range("E23").Value = "Zombie"
$wsMonsters.Range("A24").Value = "Bear"
$wsMonsters.Range("C24").Value = "Steel Golem"
$wsMonsters.Range("E24").Value = "ghost"
$wsMonsters.Range("A25").Value = "Arctic Bear"
$wsMonsters.Range("C25").Value = "Mithril Golem"
$wsMonsters.Range("E25").Value = "vampire"
$wsMonsters.Range("E26").Value = "Dullahan"

# ---------------------------------------------------------------
# 4. Create the new blank 'Sheet1' worksheet, positioned right
#    after 'Names'. Creating it second means it gets sheetId 11.
# ---------------------------------------------------------------
$wsBlank = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Names"))
$wsBlank.Name = "Sheet1"

# ---------------------------------------------------------------
# 5. Create the new 'Corruptions' worksheet, placed at the very
#    end (after 'Monsters'). It gets sheetId 12.
# ---------------------------------------------------------------
$wsCorruptions = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsCorruptions.Name = "Corruptions"
$wsCorruptions.Range("A1").Value = "Name"
$wsCorruptions.Range("D1").Value = "Effect"
$wsCorruptions.Range("K1").Value = "Description"
$wsCorruptions.Range("A3").Value = "Of Two Minds"
$wsCorruptions.Range("D3").Value = "gain an extra turn of movement that is randomly used"
$wsCorruptions.Range("K3").Value = "You feel as if your mind was split in two"
$wsCorruptions.Range("D4").Value = " -2 int, -1 learning"
$wsCorruptions.Range("A6").Value = "Of 2 Minds"
$wsCorruptions.Range("D6").Value = "gain an extra turn of movement that is randomly used"
$wsCorruptions.Range("K6").Value = "You feel as if your mind has abosrbed another entitiy"
$wsCorruptions.Range("D7").Value = " +2 int, + 1 learning"

# ---------------------------------------------------------------
# 6. Restore/update sheet selections & the active tab so the
#    saved view state matches the edited workbook.
# ---------------------------------------------------------------
$wsLinks.Activate()
$wsLinks.Range("A6").Select()

$wsCharCreation = $wb.Worksheets.Item("Character Creation")
$wsCharCreation.Activate()
$wsCharCreation.Range("H30").Select()

$wsNames.Activate()
$wsNames.Range("Q39").Select()

$wsMonsters.Activate()
$wsMonsters.Range("C9").Select()

$wsCorruptions.Activate()
$wsCorruptions.Range("G15").Select()

try {
    $win = $excel.ActiveWindow
    $win.TabRatio = 0.731
} catch {}

Write-Output ("Final sheet count: " + $wb.Worksheets.Count)
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Output ($i.ToString() + ": " + $wb.Worksheets.Item($i).Name)
}
